$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.9938993333333332
$ws.Range("H2").Value = 2.981698
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 4.807635
$ws.Range("N2").Value = 14.422905
$ws.Range("O2").Value = 0.1639819574772189
$ws.Range("P2").Value = 0.1639819574772189
$ws.Range("Q2").Value = 4.77830522141
$ws.Range("R2").Value = 43.00474699269
$ws.Range("S2").Value = 0.1639819574772189
$ws.Range("T2").Value = 0.1639819574772189

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.9938993333333332
$ws.Range("H3").Value = 2.981698
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 18.86426066666667
$ws.Range("N3").Value = 56.592782
$ws.Range("O3").Value = 0.643434534959602
$ws.Range("P3").Value = 0.6434345349596021
$ws.Range("Q3").Value = 18.74917610042622
$ws.Range("R3").Value = 168.742584903836
$ws.Range("S3").Value = 0.643434534959602
$ws.Range("T3").Value = 0.6434345349596021

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.9938993333333332
$ws.Range("H4").Value = 2.981698
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 5.646177333333333
$ws.Range("N4").Value = 16.938532
$ws.Range("O4").Value = 0.192583507563179
$ws.Range("P4").Value = 0.192583507563179
$ws.Range("Q4").Value = 5.611731887481777
$ws.Range("R4").Value = 50.50558698733599
$ws.Range("S4").Value = 0.192583507563179
$ws.Range("T4").Value = 0.192583507563179
